$wb = $excel.ActiveWorkbook

$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

$ws = $wb.Worksheets.Item(1)
$values1 = @(3906.399109145206, 0, 48353.76274462014, 0, 289724.0114301849, 9433.134471502228, 0, 2534.277928792104, 0, 0, 0, 0, 0, 2366.837285808575, 1995.50012305223)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
    $ws.Cells.Item(2, $c).Value = $values1[$c-1]
}

$ws = $wb.Worksheets.Item(2)
$values2 = @(6991.052031681918, 0, 197913.7502057619, 0, 289724.0114301849, 16452.51445364119, 0, 8194.52068131253, 0, 0, 0, 0, 0, 7541.022925285812, 6256.672442780481)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
    $ws.Cells.Item(2, $c).Value = $values2[$c-1]
}

$ws = $wb.Worksheets.Item(3)
$values3 = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 12883.76677115856, 9262.01660481554)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
    $ws.Cells.Item(2, $c).Value = $values3[$c-1]
}

$ws = $wb.Worksheets.Item(4)
$values4 = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 14041.61652360174, 9262.01660481554)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
    $ws.Cells.Item(2, $c).Value = $values4[$c-1]
}

$ws = $wb.Worksheets.Item(5)
$values5 = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16872.73121247132, 10094.37971814901)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
    $ws.Cells.Item(2, $c).Value = $values5[$c-1]
}

$ws = $wb.Worksheets.Item(6)
$values6 = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16872.73121247132, 10094.37971814901)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
    $ws.Cells.Item(2, $c).Value = $values6[$c-1]
}
